$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 169, shifting rows 169-220 down to 170-221.
$ws.Rows("169:169").Insert()

# Fill the new row 169 with the new record.
$ws.Range("A169").Value = 5
$ws.Range("B169").Value = "Macroferia Regional de Talca"
$ws.Range("C169").Value = "Maule"
$ws.Range("D169").Value = 44463
$ws.Range("E169").Value = 7
$ws.Range("F169").Value = 100112043
$ws.Range("G169").Value = "Pepino ensalada"
$ws.Range("H169").Value = "Sin especificar"
$ws.Range("I169").Value = "Primera"
$ws.Range("J169").Value = 300
$ws.Range("K169").Value = 16000
$ws.Range("L169").Value = 16000
$ws.Range("M169").Value = 16000
$ws.Range("N169").Value = "$/caja 60 unidades"
$ws.Range("O169").Value = "Región de Arica y Parinacota"
$ws.Range("P169").Value = 267
$ws.Range("Q169").Value = 60
$ws.Range("R169").Value = "Hortaliza"
